$d = $word.ActiveDocument

# --- Change 1: paragraph 2, "realiza" -> "hace" ---
$para1 = $d.Paragraphs(2).Range
$full1 = $para1.Text
$idx1 = $full1.IndexOf("y realiza una solicitud")
$wordStart1 = $para1.Start + $idx1 + 2
$wordEnd1 = $wordStart1 + 7

$restRange1 = $d.Range($wordStart1, $para1.End)
$restRange1.Bold = 1

$wordRange1 = $d.Range($wordStart1, $wordEnd1)
$wordRange1.Text = "hace"

$afterRange1 = $d.Range($wordStart1 + 4, $para1.End)
$afterRange1.Bold = 0
$wordRange1b = $d.Range($wordStart1, $wordStart1 + 4)
$wordRange1b.Bold = 0

# --- Change 2: paragraph 3, "localiza" -> "encuentra" ---
$para2 = $d.Paragraphs(3).Range
$full2 = $para2.Text
$idx2 = $full2.IndexOf("disponibilidad, localiza el libro")
$wordStart2 = $para2.Start + $idx2 + 16
$wordEnd2 = $wordStart2 + 8

$restRange2 = $d.Range($wordStart2, $para2.End)
$restRange2.Bold = 1

$wordRange2 = $d.Range($wordStart2, $wordEnd2)
$wordRange2.Text = "encuentra"

$afterRange2 = $d.Range($wordStart2 + 9, $para2.End)
$afterRange2.Bold = 0
$wordRange2b = $d.Range($wordStart2, $wordStart2 + 9)
$wordRange2b.Bold = 0
